# Generate Report for Handoff
# Updates the status/handoff-datetime for the "dbdd0596-2928-41d4-ac43-bdc64e6bd439" entry
# (row 3) from "In Translation" to "Ready for handoff" and refreshes the
# associated handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 06:21:53"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-24 06:21:49"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-24 06:21:53"
